# Rename the variable description text on the "About" sheet:
# "BVTStL Boolean Vehicle Types Subject to LCFS" -> "BVTStL BAU Vehicle Types Subject to LCFS"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("About")

$ws.Range("A1").Value = "BVTStL BAU Vehicle Types Subject to LCFS"

# Update the selected cell to match the saved view state (A2 on the About sheet)
$ws.Range("A2").Select() | Out-Null
